$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "User's Guide"
$ws.Range("A7").Value = "TivaWare Peripheral Driver Library"
$ws.Range("C7").Value = "TivaWare Peripheral Driver Library"
$ws.Range("D7").Value = "spmu298d.pdf"

